$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The "Value" column (G) holds numeric-looking data stored as TEXT (shared
# strings), not numbers. Assigning a numeric-looking string via .Value /
# .Value2 would make Excel convert the cell to a real number, which does not
# match the original authoring (all of column G is text). Instead, for each
# correction we copy the cell from an existing row that already holds the
# correct text value - Range.Copy preserves the source cell's text data type
# without touching any styles.

# corrections: target row -> source row that already holds the desired text
$corrections = [ordered]@{
    8   = 4    # "2" -> "1"
    14  = 21   # "4" -> "2"
    15  = 21   # "4" -> "2"
    30  = 25   # "10" -> "5"
    34  = 4    # "2" -> "1"
    291 = 4    # "2" -> "1"
    297 = 21   # "4" -> "2"
    313 = 25   # "10" -> "5"
    314 = 21   # "6" -> "2"
    317 = 21   # "4" -> "2"
    574 = 4    # "2" -> "1"
    580 = 21   # "4" -> "2"
    596 = 25   # "10" -> "5"
    597 = 21   # "6" -> "2"
    600 = 21   # "4" -> "2"
    857 = 24   # "6" -> "3"
    863 = 7    # "8" -> "4"
    864 = 4    # "2" -> "1"
    879 = 878  # "44" -> "22"
    880 = 22   # "21" -> "7"
    883 = 24   # "6" -> "3"
}

foreach ($targetRow in $corrections.Keys) {
    $sourceRow = $corrections[$targetRow]
    $ws.Range("G$sourceRow").Copy($ws.Range("G$targetRow"))
}
